# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 139
    4  = 350
    6  = 1844
    7  = 76
    8  = 114
    10 = 727
    11 = 305
    13 = 4400
    15 = 324
    16 = 1197
    17 = 512
    19 = 766
    21 = 398
    22 = 53
    23 = 204
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
